$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 4083590544.4
$ws.Range("P2").Value = 1140919499.37
$ws.Range("Q2").Value = 495525853.84
$ws.Range("R2").Value = 834.3795513188001
$ws.Range("S2").Value = 460625860.26
$ws.Range("T2").Value = 96.50516999
$ws.Range("U2").Value = 594596732.96
$ws.Range("V2").Value = 82.4105021087
$ws.Range("W2").Value = 3593550932.53
$ws.Range("X2").Value = 512747000.21
$ws.Range("Y2").Value = 281.9479676409
$ws.Range("Z2").Value = 28800574.99
$ws.Range("AA2").Value = -77.126149301
$ws.Range("AB2").Value = 490039611.87
$ws.Range("AC2").Value = 62.7675612144
$ws.Range("AD2").Value = 101.2967578078
$ws.Range("AE2").Value = 108.0113000079
$ws.Range("AF2").Value = 89.00982962800001
$ws.Range("AG2").Value = 87.99978581249999
